# Apply "Update Daily Report: 2026-02-17" edit
# Adds the next business day (2026-02-13, Excel serial 46066) of depository
# silver-stock movements to Daily_Data, and refreshes the derived pivot/summary
# tables on Today_Summary and Monthly_Stats to reflect the new day.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Daily_Data")
$ws2 = $wb.Worksheets.Item("Today_Summary")
$ws3 = $wb.Worksheets.Item("Monthly_Stats")

# --- Daily_Data: append 24 new detail rows (rows 122-145) for the new date ---
$ws1.Cells.Item(122,1).Value = 46066
$ws1.Cells.Item(122,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(122,2).Value = "ASAHI DEPOSITORY LLC Registered"
$ws1.Cells.Item(122,3).Value = 23953631.592
$ws1.Cells.Item(122,4).Value = 0
$ws1.Cells.Item(122,5).Value = 0
$ws1.Cells.Item(122,6).Value = 0
$ws1.Cells.Item(122,7).Value = 0
$ws1.Cells.Item(122,8).Value = 23953631.592

$ws1.Cells.Item(123,1).Value = 46066
$ws1.Cells.Item(123,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(123,2).Value = "ASAHI DEPOSITORY LLC Eligible"
$ws1.Cells.Item(123,3).Value = 2555897.608
$ws1.Cells.Item(123,4).Value = 0
$ws1.Cells.Item(123,5).Value = 458859.4
$ws1.Cells.Item(123,6).Value = -458859.4
$ws1.Cells.Item(123,7).Value = 0
$ws1.Cells.Item(123,8).Value = 2097038.208

$ws1.Cells.Item(124,1).Value = 46066
$ws1.Cells.Item(124,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(124,2).Value = "BRINK'S, INC. Registered"
$ws1.Cells.Item(124,3).Value = 16122359.646
$ws1.Cells.Item(124,4).Value = 0
$ws1.Cells.Item(124,5).Value = 0
$ws1.Cells.Item(124,6).Value = 0
$ws1.Cells.Item(124,7).Value = 0
$ws1.Cells.Item(124,8).Value = 16122359.646

$ws1.Cells.Item(125,1).Value = 46066
$ws1.Cells.Item(125,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(125,2).Value = "BRINK'S, INC. Eligible"
$ws1.Cells.Item(125,3).Value = 40640060.474
$ws1.Cells.Item(125,4).Value = 0
$ws1.Cells.Item(125,5).Value = 1052287.68
$ws1.Cells.Item(125,6).Value = -1052287.68
$ws1.Cells.Item(125,7).Value = 0
$ws1.Cells.Item(125,8).Value = 39587772.794

$ws1.Cells.Item(126,1).Value = 46066
$ws1.Cells.Item(126,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(126,2).Value = "CNT DEPOSITORY, INC. Registered"
$ws1.Cells.Item(126,3).Value = 12974598.079
$ws1.Cells.Item(126,4).Value = 0
$ws1.Cells.Item(126,5).Value = 0
$ws1.Cells.Item(126,6).Value = 0
$ws1.Cells.Item(126,7).Value = -739341.701
$ws1.Cells.Item(126,8).Value = 12235256.378

$ws1.Cells.Item(127,1).Value = 46066
$ws1.Cells.Item(127,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(127,2).Value = "CNT DEPOSITORY, INC. Eligible"
$ws1.Cells.Item(127,3).Value = 14961602.228
$ws1.Cells.Item(127,4).Value = 0
$ws1.Cells.Item(127,5).Value = 673725.54
$ws1.Cells.Item(127,6).Value = -673725.54
$ws1.Cells.Item(127,7).Value = 739341.701
$ws1.Cells.Item(127,8).Value = 15027218.389

$ws1.Cells.Item(128,1).Value = 46066
$ws1.Cells.Item(128,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(128,2).Value = "DELAWARE DEPOSITORY Registered"
$ws1.Cells.Item(128,3).Value = 1552701.933
$ws1.Cells.Item(128,4).Value = 0
$ws1.Cells.Item(128,5).Value = 0
$ws1.Cells.Item(128,6).Value = 0
$ws1.Cells.Item(128,7).Value = -5006.7
$ws1.Cells.Item(128,8).Value = 1547695.233

$ws1.Cells.Item(129,1).Value = 46066
$ws1.Cells.Item(129,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(129,2).Value = "DELAWARE DEPOSITORY Eligible"
$ws1.Cells.Item(129,3).Value = 16249560.362
$ws1.Cells.Item(129,4).Value = 0
$ws1.Cells.Item(129,5).Value = 0
$ws1.Cells.Item(129,6).Value = 0
$ws1.Cells.Item(129,7).Value = 5006.7
$ws1.Cells.Item(129,8).Value = 16254567.062

$ws1.Cells.Item(130,1).Value = 46066
$ws1.Cells.Item(130,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(130,2).Value = "HSBC BANK, USA Registered"
$ws1.Cells.Item(130,3).Value = 3472271.68
$ws1.Cells.Item(130,4).Value = 0
$ws1.Cells.Item(130,5).Value = 0
$ws1.Cells.Item(130,6).Value = 0
$ws1.Cells.Item(130,7).Value = 0
$ws1.Cells.Item(130,8).Value = 3472271.68

$ws1.Cells.Item(131,1).Value = 46066
$ws1.Cells.Item(131,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(131,2).Value = "HSBC BANK, USA Eligible"
$ws1.Cells.Item(131,3).Value = 21150312.483
$ws1.Cells.Item(131,4).Value = 0
$ws1.Cells.Item(131,5).Value = 0
$ws1.Cells.Item(131,6).Value = 0
$ws1.Cells.Item(131,7).Value = 0
$ws1.Cells.Item(131,8).Value = 21150312.483

$ws1.Cells.Item(132,1).Value = 46066
$ws1.Cells.Item(132,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(132,2).Value = "INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Registered"
$ws1.Cells.Item(132,3).Value = 273789.87
$ws1.Cells.Item(132,4).Value = 0
$ws1.Cells.Item(132,5).Value = 0
$ws1.Cells.Item(132,6).Value = 0
$ws1.Cells.Item(132,7).Value = 0
$ws1.Cells.Item(132,8).Value = 273789.87

$ws1.Cells.Item(133,1).Value = 46066
$ws1.Cells.Item(133,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(133,2).Value = "INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Eligible"
$ws1.Cells.Item(133,3).Value = 3642206.244
$ws1.Cells.Item(133,4).Value = 0
$ws1.Cells.Item(133,5).Value = 0
$ws1.Cells.Item(133,6).Value = 0
$ws1.Cells.Item(133,7).Value = 0
$ws1.Cells.Item(133,8).Value = 3642206.244

$ws1.Cells.Item(134,1).Value = 46066
$ws1.Cells.Item(134,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(134,2).Value = "JP MORGAN CHASE BANK NA Registered"
$ws1.Cells.Item(134,3).Value = 12035257.32
$ws1.Cells.Item(134,4).Value = 0
$ws1.Cells.Item(134,5).Value = 0
$ws1.Cells.Item(134,6).Value = 0
$ws1.Cells.Item(134,7).Value = -10008.78
$ws1.Cells.Item(134,8).Value = 12025248.54

$ws1.Cells.Item(135,1).Value = 46066
$ws1.Cells.Item(135,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(135,2).Value = "JP MORGAN CHASE BANK NA Eligible"
$ws1.Cells.Item(135,3).Value = 147710753.583
$ws1.Cells.Item(135,4).Value = 0
$ws1.Cells.Item(135,5).Value = 1947445.7
$ws1.Cells.Item(135,6).Value = -1947445.7
$ws1.Cells.Item(135,7).Value = 10008.78
$ws1.Cells.Item(135,8).Value = 145773316.663

$ws1.Cells.Item(136,1).Value = 46066
$ws1.Cells.Item(136,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(136,2).Value = "LOOMIS INTERNATIONAL (US) LLC Registered"
$ws1.Cells.Item(136,3).Value = 7374299.767
$ws1.Cells.Item(136,4).Value = 0
$ws1.Cells.Item(136,5).Value = 0
$ws1.Cells.Item(136,6).Value = 0
$ws1.Cells.Item(136,7).Value = 19054.65
$ws1.Cells.Item(136,8).Value = 7393354.417

$ws1.Cells.Item(137,1).Value = 46066
$ws1.Cells.Item(137,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(137,2).Value = "LOOMIS INTERNATIONAL (US) LLC Eligible"
$ws1.Cells.Item(137,3).Value = 23295383.436
$ws1.Cells.Item(137,4).Value = 0
$ws1.Cells.Item(137,5).Value = 324212.08
$ws1.Cells.Item(137,6).Value = -324212.08
$ws1.Cells.Item(137,7).Value = -19054.65
$ws1.Cells.Item(137,8).Value = 22952116.706

$ws1.Cells.Item(138,1).Value = 46066
$ws1.Cells.Item(138,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(138,2).Value = "MALCA-AMIT ARMORED, INC. Registered"
$ws1.Cells.Item(138,3).Value = 0
$ws1.Cells.Item(138,4).Value = 0
$ws1.Cells.Item(138,5).Value = 0
$ws1.Cells.Item(138,6).Value = 0
$ws1.Cells.Item(138,7).Value = 0
$ws1.Cells.Item(138,8).Value = 0

$ws1.Cells.Item(139,1).Value = 46066
$ws1.Cells.Item(139,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(139,2).Value = "MALCA-AMIT ARMORED, INC. Eligible"
$ws1.Cells.Item(139,3).Value = 0
$ws1.Cells.Item(139,4).Value = 0
$ws1.Cells.Item(139,5).Value = 0
$ws1.Cells.Item(139,6).Value = 0
$ws1.Cells.Item(139,7).Value = 0
$ws1.Cells.Item(139,8).Value = 0

$ws1.Cells.Item(140,1).Value = 46066
$ws1.Cells.Item(140,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(140,2).Value = "MALCA-AMIT USA, LLC Registered"
$ws1.Cells.Item(140,3).Value = 1225506.264
$ws1.Cells.Item(140,4).Value = 0
$ws1.Cells.Item(140,5).Value = 0
$ws1.Cells.Item(140,6).Value = 0
$ws1.Cells.Item(140,7).Value = 0
$ws1.Cells.Item(140,8).Value = 1225506.264

$ws1.Cells.Item(141,1).Value = 46066
$ws1.Cells.Item(141,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(141,2).Value = "MALCA-AMIT USA, LLC Eligible"
$ws1.Cells.Item(141,3).Value = 798026.177
$ws1.Cells.Item(141,4).Value = 0
$ws1.Cells.Item(141,5).Value = 0
$ws1.Cells.Item(141,6).Value = 0
$ws1.Cells.Item(141,7).Value = 0
$ws1.Cells.Item(141,8).Value = 798026.177

$ws1.Cells.Item(142,1).Value = 46066
$ws1.Cells.Item(142,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(142,2).Value = "MANFRA, TORDELLA & BROOKES, LLC Registered"
$ws1.Cells.Item(142,3).Value = 6370259.821
$ws1.Cells.Item(142,4).Value = 0
$ws1.Cells.Item(142,5).Value = 0
$ws1.Cells.Item(142,6).Value = 0
$ws1.Cells.Item(142,7).Value = -4827.529
$ws1.Cells.Item(142,8).Value = 6365432.292

$ws1.Cells.Item(143,1).Value = 46066
$ws1.Cells.Item(143,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(143,2).Value = "MANFRA, TORDELLA & BROOKES, LLC Eligible"
$ws1.Cells.Item(143,3).Value = 12298021.519
$ws1.Cells.Item(143,4).Value = 0
$ws1.Cells.Item(143,5).Value = 0
$ws1.Cells.Item(143,6).Value = 0
$ws1.Cells.Item(143,7).Value = 4827.529
$ws1.Cells.Item(143,8).Value = 12302849.048

$ws1.Cells.Item(144,1).Value = 46066
$ws1.Cells.Item(144,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(144,2).Value = "STONEX PRECIOUS METALS LLC Registered"
$ws1.Cells.Item(144,3).Value = 7545291.14
$ws1.Cells.Item(144,4).Value = 0
$ws1.Cells.Item(144,5).Value = 0
$ws1.Cells.Item(144,6).Value = 0
$ws1.Cells.Item(144,7).Value = -4967.6
$ws1.Cells.Item(144,8).Value = 7540323.54

$ws1.Cells.Item(145,1).Value = 46066
$ws1.Cells.Item(145,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(145,2).Value = "STONEX PRECIOUS METALS LLC Eligible"
$ws1.Cells.Item(145,3).Value = 233197.38
$ws1.Cells.Item(145,4).Value = 0
$ws1.Cells.Item(145,5).Value = 4967.8
$ws1.Cells.Item(145,6).Value = -4967.8
$ws1.Cells.Item(145,7).Value = 4967.6
$ws1.Cells.Item(145,8).Value = 233197.18

# --- Today_Summary: refresh per-depository Eligible / Registered / Total_Stock ---
$ws2.Cells.Item(2,2).Value = 2097038.208
$ws2.Cells.Item(2,3).Value = 23953631.592
$ws2.Cells.Item(2,4).Value = 26050669.8
$ws2.Cells.Item(3,2).Value = 39587772.794
$ws2.Cells.Item(3,3).Value = 16122359.646
$ws2.Cells.Item(3,4).Value = 55710132.44
$ws2.Cells.Item(4,2).Value = 15027218.389
$ws2.Cells.Item(4,3).Value = 12235256.378
$ws2.Cells.Item(4,4).Value = 27262474.767
$ws2.Cells.Item(5,2).Value = 16254567.062
$ws2.Cells.Item(5,3).Value = 1547695.233
$ws2.Cells.Item(5,4).Value = 17802262.295
$ws2.Cells.Item(6,2).Value = 21150312.483
$ws2.Cells.Item(6,3).Value = 3472271.68
$ws2.Cells.Item(6,4).Value = 24622584.163
$ws2.Cells.Item(7,2).Value = 3642206.244
$ws2.Cells.Item(7,3).Value = 273789.87
$ws2.Cells.Item(7,4).Value = 3915996.114
$ws2.Cells.Item(8,2).Value = 145773316.663
$ws2.Cells.Item(8,3).Value = 12025248.54
$ws2.Cells.Item(8,4).Value = 157798565.203
$ws2.Cells.Item(9,2).Value = 22952116.706
$ws2.Cells.Item(9,3).Value = 7393354.417
$ws2.Cells.Item(9,4).Value = 30345471.123
$ws2.Cells.Item(10,2).Value = 0
$ws2.Cells.Item(10,3).Value = 0
$ws2.Cells.Item(10,4).Value = 0
$ws2.Cells.Item(11,2).Value = 798026.177
$ws2.Cells.Item(11,3).Value = 1225506.264
$ws2.Cells.Item(11,4).Value = 2023532.441
$ws2.Cells.Item(12,2).Value = 12302849.048
$ws2.Cells.Item(12,3).Value = 6365432.292
$ws2.Cells.Item(12,4).Value = 18668281.34
$ws2.Cells.Item(13,2).Value = 233197.18
$ws2.Cells.Item(13,3).Value = 7540323.54
$ws2.Cells.Item(13,4).Value = 7773520.72

# --- Monthly_Stats: refresh month-to-date Eligible / Registered / Grand_Total ---
$ws3.Cells.Item(2,2).Value = 279818620.954
$ws3.Cells.Item(2,3).Value = 92154869.452
$ws3.Cells.Item(2,4).Value = 371973490.406

# --- Monthly_Stats: refresh per-depository/type cumulative WITHDRAWN and latest TOTAL_TODAY ---
$ws3.Cells.Item(7,4).Value = 2130390.53
$ws3.Cells.Item(7,5).Value = 2097038.208
$ws3.Cells.Item(8,4).Value = 0
$ws3.Cells.Item(8,5).Value = 23953631.592
$ws3.Cells.Item(9,4).Value = 1866753.425
$ws3.Cells.Item(9,5).Value = 39587772.794
$ws3.Cells.Item(10,4).Value = 0
$ws3.Cells.Item(10,5).Value = 16122359.646
$ws3.Cells.Item(11,4).Value = 3915027.358
$ws3.Cells.Item(11,5).Value = 15027218.389
$ws3.Cells.Item(12,4).Value = 0
$ws3.Cells.Item(12,5).Value = 12235256.378
$ws3.Cells.Item(13,4).Value = 186666.419
$ws3.Cells.Item(13,5).Value = 16254567.062
$ws3.Cells.Item(14,4).Value = 0
$ws3.Cells.Item(14,5).Value = 1547695.233
$ws3.Cells.Item(15,4).Value = 110629.57
$ws3.Cells.Item(15,5).Value = 21150312.483
$ws3.Cells.Item(16,4).Value = 0
$ws3.Cells.Item(16,5).Value = 3472271.68
$ws3.Cells.Item(17,4).Value = 0
$ws3.Cells.Item(17,5).Value = 3642206.244
$ws3.Cells.Item(18,4).Value = 0
$ws3.Cells.Item(18,5).Value = 273789.87
$ws3.Cells.Item(19,4).Value = 11065580
$ws3.Cells.Item(19,5).Value = 145773316.663
$ws3.Cells.Item(20,4).Value = 0
$ws3.Cells.Item(20,5).Value = 12025248.54
$ws3.Cells.Item(21,4).Value = 2964461.58
$ws3.Cells.Item(21,5).Value = 22952116.706
$ws3.Cells.Item(22,4).Value = 0
$ws3.Cells.Item(22,5).Value = 7393354.417
$ws3.Cells.Item(23,4).Value = 0
$ws3.Cells.Item(23,5).Value = 0
$ws3.Cells.Item(24,4).Value = 0
$ws3.Cells.Item(24,5).Value = 0
$ws3.Cells.Item(25,4).Value = 0
$ws3.Cells.Item(25,5).Value = 798026.177
$ws3.Cells.Item(26,4).Value = 0
$ws3.Cells.Item(26,5).Value = 1225506.264
$ws3.Cells.Item(27,4).Value = 849925.823
$ws3.Cells.Item(27,5).Value = 12302849.048
$ws3.Cells.Item(28,4).Value = 0
$ws3.Cells.Item(28,5).Value = 6365432.292
$ws3.Cells.Item(29,4).Value = 9949.82
$ws3.Cells.Item(29,5).Value = 233197.18
$ws3.Cells.Item(30,4).Value = 0
$ws3.Cells.Item(30,5).Value = 7540323.54
